$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row data (Element / ElementType / Locator) ---
$ws.Cells.Item(1,1).Value = 'Element'
$ws.Cells.Item(1,2).Value = 'ElementType'
$ws.Cells.Item(1,3).Value = 'Locator'

$ws.Cells.Item(2,1).Value = 'userName'
$ws.Cells.Item(2,2).Value = 'textbox'
$ws.Cells.Item(2,3).Value = "xp&//input[@name='email']"

$ws.Cells.Item(3,1).Value = 'password'
$ws.Cells.Item(3,2).Value = 'textbox'
$ws.Cells.Item(3,3).Value = "xp&//input[@name='password']"

$ws.Cells.Item(4,1).Value = 'loginbutton'
$ws.Cells.Item(4,2).Value = 'button'
$ws.Cells.Item(4,3).Value = 'bt&Sign in'

$ws.Cells.Item(5,1).Value = 'validateloginmessage'
$ws.Cells.Item(5,2).Value = 'validatetext'
$ws.Cells.Item(5,3).Value = 'div.alert-message'

$ws.Cells.Item(6,1).Value = 'AcceptLoginAlert'
$ws.Cells.Item(6,2).Value = 'button'
$ws.Cells.Item(6,3).Value = 'css&button.alert-button'

$ws.Cells.Item(7,1).Value = 'LoanPurpose'
$ws.Cells.Item(7,2).Value = 'button'
$ws.Cells.Item(7,3).Value = 'css&ion-segment[name="purchaseOrRefi"] ion-segment-button:nth-child(2)'

$ws.Cells.Item(8,1).Value = 'PhoneNumber'
$ws.Cells.Item(8,2).Value = 'textbox'
$ws.Cells.Item(8,3).Value = "xp&//input[@name='phoneNumber']"

$ws.Cells.Item(9,1).Value = 'OPT_contact Policy'
$ws.Cells.Item(9,2).Value = 'button'
$ws.Cells.Item(9,3).Value = "xp&//ion-segment-button[text()='No']"

$ws.Cells.Item(10,1).Value = 'SaveLoanPurpose'
$ws.Cells.Item(10,2).Value = 'button'
$ws.Cells.Item(10,3).Value = 'css&button.app-button.save-continue'

$ws.Cells.Item(11,1).Value = 'validateContactMessage'
$ws.Cells.Item(11,2).Value = 'validatetext'
$ws.Cells.Item(11,3).Value = 'css&div.error-box ion-content ion-label'

$ws.Cells.Item(12,1).Value = 'acceptErrors'
$ws.Cells.Item(12,2).Value = 'button'
$ws.Cells.Item(12,3).Value = 'css&bl-validation-error div.error-box div.page-bar button'

$ws.Cells.Item(13,1).Value = 'openAdressForm'
$ws.Cells.Item(13,2).Value = 'button'
$ws.Cells.Item(13,3).Value = "xp&//ion-input[@name='address']"

$ws.Cells.Item(14,1).Value = 'Adress Line'
$ws.Cells.Item(14,2).Value = 'textbox'
$ws.Cells.Item(14,3).Value = 'css&input[name="addressLine"]'

$ws.Cells.Item(15,1).Value = 'State'
$ws.Cells.Item(15,2).Value = 'textbox'
$ws.Cells.Item(15,3).Value = 'css&input[name="city"]'

$ws.Cells.Item(16,1).Value = 'City '
$ws.Cells.Item(16,2).Value = 'textbox'
$ws.Cells.Item(16,3).Value = 'css&input[name="state"]'

$ws.Cells.Item(17,1).Value = 'Zip'
$ws.Cells.Item(17,2).Value = 'textbox'
$ws.Cells.Item(17,3).Value = 'css&input[name="zip"]'

$ws.Cells.Item(18,1).Value = 'closeAdressForm'
$ws.Cells.Item(18,2).Value = 'button'
$ws.Cells.Item(18,3).Value = 'css&button.auto-done'

$ws.Cells.Item(19,1).Value = 'PrimaryResidence'
$ws.Cells.Item(19,2).Value = 'button'
$ws.Cells.Item(19,3).Value = 'css&ion-toggle[name="primaryResidence"]'

# Row 5 column A ("validateloginmessage") picks up the same bold + text
# number-format style already used by A2 ("userName").
$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Font.Bold = $true

# --- Column widths / best-fit ---
$ws.Columns.Item(1).ColumnWidth = 22.1
$ws.Columns.Item(3).ColumnWidth = 68.8

# --- View / selection ---
[void]$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("D2").Select()
